$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.226.60"
$ws.Range("E2").Value = "'  +0.80%  "
$ws.Range("D3").Value = "'3.471.69"
$ws.Range("E3").Value = "'  -0.40%  "
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("D5").Value = "'592.74"
$ws.Range("E5").Value = "'  +0.17%  "
$ws.Range("D6").Value = "'177.96"
$ws.Range("E6").Value = "'  +3.84%  "
$ws.Range("D8").Value = "'3.471.97"
$ws.Range("E8").Value = "'  -0.38%  "
$ws.Range("D9").Value = "'0.590"
$ws.Range("E9").Value = "'  -0.24%  "
$ws.Range("E10").Value = "'  +5.20%  "
$ws.Range("E11").Value = "'  -2.39%  "
$ws.Range("E12").Value = "'  +0.25%  "
$ws.Range("D13").Value = "'4.074.39"
$ws.Range("E13").Value = "'  -0.36%  "
$ws.Range("D14").Value = "'31.88"
$ws.Range("E14").Value = "'  +10.68%  "
$ws.Range("E15").Value = "'  +1.50%  "
$ws.Range("D16").Value = "'67.278.01"
$ws.Range("E16").Value = "'  +0.82%  "
$ws.Range("E17").Value = "'  -0.19%  "
$ws.Range("D18").Value = "'3.477.20"
$ws.Range("E18").Value = "'  -0.11%  "
$ws.Range("E19").Value = "'  -0.34%  "
$ws.Range("D20").Value = "'14.23"
$ws.Range("E20").Value = "'  +1.38%  "
$ws.Range("D21").Value = "'388.26"
$ws.Range("E21").Value = "'  -0.84%  "
$ws.Range("D22").Value = "'7.85"
$ws.Range("E22").Value = "'  -0.56%  "
$ws.Range("D23").Value = "'73.44"
$ws.Range("E23").Value = "'  +1.12%  "
$ws.Range("D24").Value = "'0.997"
$ws.Range("E24").Value = "'  -0.28%  "
$ws.Range("D25").Value = "'5.72"
$ws.Range("E25").Value = "'  +0.67%  "
$ws.Range("D26").Value = "'0.534"
$ws.Range("E26").Value = "'  -0.02%  "
$ws.Range("E27").Value = "'  +0.62%  "
$ws.Range("D28").Value = "'10.32"
$ws.Range("E29").Value = "'  -3.10%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "'  +0.16%  "
$ws.Range("E31").Value = "'  -0.61%  "
$ws.Range("D32").Value = "'1.41"
$ws.Range("E32").Value = "'  -0.23%  "
$ws.Range("E33").Value = "'  +0.37%  "
$ws.Range("E34").Value = "'  -0.72%  "
$ws.Range("D35").Value = "'7.34"
$ws.Range("E35").Value = "'  +0.59%  "
$ws.Range("E36").Value = "'  +0.04%  "
$ws.Range("D37").Value = "'1.58"
$ws.Range("E37").Value = "'  -1.80%  "
$ws.Range("D38").Value = "'163.72"
$ws.Range("E38").Value = "'  +0.59%  "
$ws.Range("D39").Value = "'0.869"
$ws.Range("E39").Value = "'  -0.75%  "
$ws.Range("E40").Value = "'  -0.35%  "
$ws.Range("D41").Value = "'2.72"
$ws.Range("E41").Value = "'  +6.50%  "
$ws.Range("D42").Value = "'6.83"
$ws.Range("E42").Value = "'  -0.57%  "
$ws.Range("D43").Value = "'4.61"
$ws.Range("E43").Value = "'  -0.72%  "
$ws.Range("D44").Value = "'26.31"
$ws.Range("E44").Value = "'  +1.20%  "
$ws.Range("D45").Value = "'2.817.25"
$ws.Range("E45").Value = "'  +0.93%  "
$ws.Range("B46").Value = "'InjectiveProtocol"
$ws.Range("C46").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'26.72"
$ws.Range("E46").Value = "'  -1.67%  "
$ws.Range("B47").Value = "'Hedera"
$ws.Range("C47").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "'0.0718"
$ws.Range("E47").Value = "'  -2.63%  "
$ws.Range("D48").Value = "'41.47"
$ws.Range("E48").Value = "'  -2.76%  "
$ws.Range("E49").Value = "'  -1.14%  "
$ws.Range("D50").Value = "'336.18"
$ws.Range("E50").Value = "'  +0.22%  "
$ws.Range("E51").Value = "'  -2.31%  "
